$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Range("D$r").Value = 44245.50784877911
}
for ($r = 16; $r -le 29; $r++) {
    $ws.Range("D$r").Value = 44245.48657256945
}
for ($r = 30; $r -le 43; $r++) {
    $ws.Range("D$r").Value = 44245.46528282407
}
